$wb = $excel.ActiveWorkbook

# --- Sheet "inputdisp" (sheet1) : tech catalog "ID" ---
$wsInput = $wb.Worksheets.Item("inputdisp")
$wsInput.Activate()

$wsInput.Range("C2").Value = 0.4
$wsInput.Range("G2").Value = 3
$wsInput.Range("H2").Value = 0.6
$wsInput.Range("I2").Value = "t2"

$wsInput.Range("G3").Value = 6
$wsInput.Range("H3").Value = 0.6
$wsInput.Range("I3").Value = "t1"

$wsInput.Range("G4").Select()

# --- Sheet "endofpipe" (sheet2) : tech catalog "EOP" ---
$wsEop = $wb.Worksheets.Item("endofpipe")
$wsEop.Activate()

$wsEop.Range("E2").Value = 0.5
$wsEop.Range("G2").Value = 13

$wsEop.Range("E3").Value = 0.2
$wsEop.Range("G3").Value = 16

$wsEop.Range("G3").Select()

# Restore the originally active sheet/tab.
$wsInput.Activate()
